# Melhoria do fluxo alternativo 7
# Reorders the customer-type step text/result across TC1-TC5 and removes
# the now-redundant trailing test-case block (rows 94-98).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TC1 (rows 16-17): was "altera para tipo de cliente B" / qty < 100
$ws.Range("B16").Value = "Usuário do Sistema mantém seleção padrão do tipo de cliente A"
$ws.Range("B17").Value = "Usuário do Sistema informa a quantidade de produtos"
$ws.Range("D17").Value = "SYSTEM registra a quantidade informada"

# TC2 (rows 35-36): was "mantém seleção padrão do tipo de cliente A" / qty >= 1000 (unchanged)
$ws.Range("B35").Value = "Usuário do Sistema seleciona tipo de cliente se desejar alterar"

# TC3 (rows 54-55): was "altera para tipo de cliente C" / qty generic
$ws.Range("B54").Value = "Usuário do Sistema altera para tipo de cliente B"
$ws.Range("B55").Value = "Usuário do Sistema informa a quantidade de produtos entre 100 e 999 unidades"
$ws.Range("D55").Value = "SYSTEM aplica fator de desconto para 100 <= quantidade < 1000: Cliente A (0,95), B (0,90), C (0,85)"

# TC4 (rows 73-74): was "seleciona tipo de cliente se desejar alterar" / qty <= 0
$ws.Range("B73").Value = "Usuário do Sistema altera para tipo de cliente C"
$ws.Range("B74").Value = "Usuário do Sistema informa a quantidade de produtos menor que 100 unidades"
$ws.Range("D74").Value = "SYSTEM aplica fator de desconto para quantidade < 100: Cliente A (0,90), B (0,85), C (0,80)"

# TC5 (rows 92-93): was "altera para tipo de cliente B" / qty 100-999
$ws.Range("B92").Value = "Usuário do Sistema mantém seleção padrão do tipo de cliente A"
$ws.Range("B93").Value = "Usuário do Sistema informa a quantidade de produtos menor ou igual a zero"
$ws.Range("D93").Value = "SYSTEM exibe mensagem 'A quantidade informada deve ser maior ou igual a 01 (um)!' (MSG002)"

# Delete the now-duplicated trailing rows (old TC5 steps 9-13 + postcondition),
# which shifts the dimension from A1:F98 to A1:F93 and removes the B98:F98 merge.
$rng = $ws.Range("A94:F98")
$rng.EntireRow.Delete()
